# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 12 - Alemania
$ws.Range("D12").Value = 171600
$ws.Range("E12").Value = 6344

# Row 56 - Kazajistan
$ws.Range("B56").Value = 13872
$ws.Range("C56").Value = 314
$ws.Range("E56").Value = 5212

# Row 87 - El Salvador
$ws.Range("E87").Value = 1801
$ws.Range("G87").Value = 4
$ws.Range("H87").Value = 68

# Row 99 - Kirguistan
$ws.Range("B99").Value = 2166
$ws.Range("C99").Value = 37
$ws.Range("D99").Value = 1668
$ws.Range("E99").Value = 472

# Row 141 - Jamaica
$ws.Range("B141").Value = 611
$ws.Range("C141").Value = 6
$ws.Range("D141").Value = 408
$ws.Range("E141").Value = 193

# Row 162 - Mongolia
$ws.Range("B162").Value = 197
$ws.Range("C162").Value = 3
$ws.Range("D162").Value = 95
$ws.Range("E162").Value = 102

# Row 183 - Butan
$ws.Range("D183").Value = 19
$ws.Range("E183").Value = 43
